$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02355713883329483
$ws.Range("H2").Value = -51.15729736824031
$ws.Range("I2").Value = -5.388595340581579

$ws.Range("G3").Value = 0.0562020186959262
$ws.Range("H3").Value = 46.52501423773307

$ws.Range("G4").Value = -0.4506179472549399
$ws.Range("H4").Value = 0.4748659954656105

$ws.Range("G5").Value = -0.4537369159229148
$ws.Range("H5").Value = 5.241484818969258

$ws.Range("G6").Value = 0.2364047493654702
$ws.Range("H6").Value = 1.186953283373938

$ws.Range("G7").Value = 0.2468848717739781
$ws.Range("H7").Value = 11.92805462000246

$ws.Range("G8").Value = 0.1631718896421494
$ws.Range("H8").Value = -2.178901828741076

$ws.Range("G9").Value = 0.1791090745465888
$ws.Range("H9").Value = 4.126324005474364

$ws.Range("G10").Value = -0.006787596473818972
$ws.Range("H10").Value = -43.21389331678113

$ws.Range("G11").Value = 0.002798499789436273
$ws.Range("H11").Value = 119.0954850175792

$ws.Range("G12").Value = 0.1333231142898378
$ws.Range("H12").Value = -2.486179146830544

$ws.Range("G13").Value = 0.1474368694660218
$ws.Range("H13").Value = 18.29246598438195

$ws.Range("G14").Value = 0.246215433374287
$ws.Range("H14").Value = -0.4520383784714342

$ws.Range("G15").Value = 0.2753852792815891
$ws.Range("H15").Value = 8.987073788124961

$ws.Range("G16").Value = 0.1450696377579788
$ws.Range("H16").Value = -5.473169409324513

$ws.Range("G17").Value = 0.1416731963598719
$ws.Range("H17").Value = -6.186133675665856

$ws.Range("G18").Value = -0.007998361608780343
$ws.Range("H18").Value = 51.15476397696379

$ws.Range("G19").Value = -0.003560084408674118
$ws.Range("H19").Value = -322.7807834197287

$ws.Range("G20").Value = 0.1382315836625981
$ws.Range("H20").Value = -0.3105311414623032

$ws.Range("G21").Value = 0.1457437049806174
$ws.Range("H21").Value = 1.851385097495627

$ws.Range("G22").Value = 0.1794181005538811
$ws.Range("H22").Value = -3.65126241590734

$ws.Range("G23").Value = 0.1757700528451485
$ws.Range("H23").Value = -2.05641424984367

$ws.Range("G24").Value = -0.09728378513350405
$ws.Range("H24").Value = -3.057845498585382

$ws.Range("G25").Value = -0.092898797592864
$ws.Range("H25").Value = 6.736031599843242

$ws.Range("G26").Value = 0.2273097649192087
$ws.Range("H26").Value = -1.229978386146269

$ws.Range("G27").Value = 0.2354548525820899
$ws.Range("H27").Value = 1.239449027691362

$ws.Range("G28").Value = 0.06583917844091934
$ws.Range("H28").Value = 11.96517921762839

$ws.Range("G29").Value = 0.07379557224227279
$ws.Range("H29").Value = 4.549228028717085
